# Fill in the Week 6 activity log for Rick (Richard Dobson).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header: name + week
$ws.Range("B2").Value = "Richard Dobson"
$ws.Range("F2").Value = "Week 6"

# Monday 2/9/2019 - group work, 9am-2pm, 5 group hours
$ws.Range("A4").Value = "Work on first iteration"
$ws.Range("C4").Value = "G"
$ws.Range("D4").Value = 43710
$ws.Range("E4").Value = 0.375
$ws.Range("F4").Value = 0.58333333333333337
$ws.Range("G4").Value = 5

# Tuesday 3/9/2019 - group work, 9am-2pm, 5 group hours
$ws.Range("A5").Value = "Work on first iteration"
$ws.Range("C5").Value = "G"
$ws.Range("D5").Value = 43711
$ws.Range("E5").Value = 0.375
$ws.Range("F5").Value = 0.58333333333333337
$ws.Range("G5").Value = 5

# Thursday 5/9/2019 - individual work, 9am-2pm, 5 individual hours
$ws.Range("A6").Value = "Work on first iteration"
$ws.Range("C6").Value = "I"
$ws.Range("D6").Value = 43713
$ws.Range("E6").Value = 0.375
$ws.Range("F6").Value = 0.58333333333333337
$ws.Range("H6").Value = 5

# Friday 6/9/2019 - group work, 9am-2pm, 5 group hours
$ws.Range("A7").Value = "Work on first iteration"
$ws.Range("C7").Value = "G"
$ws.Range("D7").Value = 43714
$ws.Range("E7").Value = 0.375
$ws.Range("F7").Value = 0.58333333333333337
$ws.Range("G7").Value = 5

# Leave cursor on the week cell, like the saved workbook.
$null = $ws.Range("F2").Select()
